$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = 1.5
$ws.Range("D13").Value = "Refactorizare capitolul 2, 3"
$ws.Range("G13").Select()
